$d = $word.ActiveDocument

# 1. Update the date: 05 de junio de 2017. -> 28 de junio de 2017.
$d.Content.Find.Execute("05 de junio de 2017.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "28 de junio de 2017.", 2)

# 2. Update the Fiscalia dependency text: Nro 3 -> Nro 1
$d.Content.Find.Execute("Fiscalia Nacional en lo Criminal y Correccional Nro 3", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Fiscalia Nacional en lo Criminal y Correccional Nro 1", 2)

# 3. Update the requested amount: $ 500,00 -> $ 942,00
$d.Content.Find.Execute("$ 500,00", $true, $false, $false, $false, $false,
                         $true, 1, $false, "$ 942,00", 2)

# 4. Fill in the justification paragraph (currently a single blank space)
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "" -and $p.Range.Text.Length -eq 2) {
        $r = $p.Range
        [void]$r.MoveEnd(1, -1)
        $r.Text = "Finalmente, la presente erogación de fondos es solicitada por este curso debido a que Es un pedido urgente"
    }
}
